$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 298.25  # H2
$ws.Cells.Item(2, 9).Value = 298.25  # I2
$ws.Cells.Item(2, 11).Value = 298.25  # K2
$ws.Cells.Item(2, 13).Value = -185.25  # M2
$ws.Cells.Item(17, 8).Value = 1368.0834  # H17
$ws.Cells.Item(17, 10).Value = 1368.0834  # J17
$ws.Cells.Item(17, 12).Value = 4104.2502  # L17
$ws.Cells.Item(17, 14).Value = -4440.2502  # N17
$ws.Cells.Item(62, 8).Value = 2550  # H62
$ws.Cells.Item(62, 9).Value = 2100  # I62
$ws.Cells.Item(62, 10).Value = 3000  # J62
$ws.Cells.Item(62, 11).Value = 2100  # K62
$ws.Cells.Item(62, 12).Value = 3000  # L62
$ws.Cells.Item(62, 13).Value = -1476  # M62
$ws.Cells.Item(62, 14).Value = -4248  # N62
$ws.Cells.Item(65, 8).Value = 2550  # H65
$ws.Cells.Item(65, 9).Value = 2100  # I65
$ws.Cells.Item(65, 10).Value = 3000  # J65
$ws.Cells.Item(65, 11).Value = 10500  # K65
$ws.Cells.Item(65, 12).Value = 15000  # L65
$ws.Cells.Item(65, 13).Value = -7380  # M65
$ws.Cells.Item(65, 14).Value = -21240  # N65
$ws.Cells.Item(135, 8).Value = 861.3333  # H135
$ws.Cells.Item(135, 9).Value = 833.6  # I135
$ws.Cells.Item(135, 10).Value = 1000  # J135
$ws.Cells.Item(135, 11).Value = 7502.400000000001  # K135
$ws.Cells.Item(135, 12).Value = 9000  # L135
$ws.Cells.Item(135, 13).Value = -4967.400000000001  # M135
$ws.Cells.Item(135, 14).Value = -14070  # N135
$ws.Cells.Item(137, 8).Value = 3537.5715  # H137
$ws.Cells.Item(137, 9).Value = 2594.6667  # I137
$ws.Cells.Item(137, 10).Value = 4244.75  # J137
$ws.Cells.Item(137, 11).Value = 7784.000100000001  # K137
$ws.Cells.Item(137, 12).Value = 12734.25  # L137
$ws.Cells.Item(137, 13).Value = -5234.000100000001  # M137
$ws.Cells.Item(137, 14).Value = -17834.25  # N137

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1383.3334  # H2
$ws.Cells.Item(2, 9).Value = 1387.2142  # I2
$ws.Cells.Item(2, 11).Value = 1387.2142  # K2
$ws.Cells.Item(2, 13).Value = -1274.2142  # M2
$ws.Cells.Item(63, 8).Value = 15290  # H63
$ws.Cells.Item(63, 9).Value = 0  # I63
$ws.Cells.Item(63, 10).Value = 15290  # J63
$ws.Cells.Item(63, 11).Value = 0  # K63
$ws.Cells.Item(63, 12).Value = 15290  # L63
$ws.Cells.Item(63, 13).ClearContents()  # M63
$ws.Cells.Item(63, 14).Value = -16662  # N63
$ws.Cells.Item(66, 8).Value = 15290  # H66
$ws.Cells.Item(66, 9).Value = 0  # I66
$ws.Cells.Item(66, 10).Value = 15290  # J66
$ws.Cells.Item(66, 11).Value = 0  # K66
$ws.Cells.Item(66, 12).Value = 76450  # L66
$ws.Cells.Item(66, 13).ClearContents()  # M66
$ws.Cells.Item(66, 14).Value = -83314  # N66
$ws.Cells.Item(74, 8).Value = 3508.25  # H74
$ws.Cells.Item(74, 9).Value = 3011  # I74
$ws.Cells.Item(74, 11).Value = 3011  # K74
$ws.Cells.Item(74, 13).Value = -2137  # M74
$ws.Cells.Item(77, 8).Value = 3508.25  # H77
$ws.Cells.Item(77, 9).Value = 3011  # I77
$ws.Cells.Item(77, 11).Value = 15055  # K77
$ws.Cells.Item(77, 13).Value = -10687  # M77
$ws.Cells.Item(116, 8).Value = 1383.3334  # H116
$ws.Cells.Item(116, 9).Value = 1387.2142  # I116
$ws.Cells.Item(116, 11).Value = 1387.2142  # K116
$ws.Cells.Item(116, 13).Value = 906.7858000000001  # M116
$ws.Cells.Item(132, 8).Value = 4251  # H132
$ws.Cells.Item(132, 9).Value = 4004  # I132
$ws.Cells.Item(132, 11).Value = 12012  # K132
$ws.Cells.Item(132, 13).Value = -9482  # M132

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1383.3334  # H3
$ws.Cells.Item(3, 9).Value = 1387.2142  # I3
$ws.Cells.Item(3, 11).Value = 1387.2142  # K3
$ws.Cells.Item(3, 13).Value = -1273.2142  # M3

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5249.75  # H31
$ws.Cells.Item(31, 9).Value = 0  # I31
$ws.Cells.Item(31, 11).Value = 0  # K31
$ws.Cells.Item(31, 13).ClearContents()  # M31
$ws.Cells.Item(34, 8).Value = 5249.75  # H34
$ws.Cells.Item(34, 9).Value = 0  # I34
$ws.Cells.Item(34, 11).Value = 0  # K34
$ws.Cells.Item(34, 13).ClearContents()  # M34
$ws.Cells.Item(58, 8).Value = 4226.857  # H58
$ws.Cells.Item(58, 9).Value = 4098  # I58
$ws.Cells.Item(58, 10).Value = 5000  # J58
$ws.Cells.Item(58, 11).Value = 4098  # K58
$ws.Cells.Item(58, 12).Value = 5000  # L58
$ws.Cells.Item(58, 13).Value = -3895  # M58
$ws.Cells.Item(58, 14).Value = -5406  # N58
$ws.Cells.Item(94, 8).Value = 4605.2  # H94
$ws.Cells.Item(94, 10).Value = 10014  # J94
$ws.Cells.Item(94, 12).Value = 10014  # L94
$ws.Cells.Item(94, 14).Value = -10916  # N94
$ws.Cells.Item(122, 8).Value = 1375.1666  # H122
$ws.Cells.Item(122, 9).Value = 1375.1666  # I122
$ws.Cells.Item(122, 11).Value = 4125.4998  # K122
$ws.Cells.Item(122, 13).Value = -1675.4998  # M122
$ws.Cells.Item(132, 8).Value = 146380.14  # H132
$ws.Cells.Item(132, 9).Value = 202532.6  # I132
$ws.Cells.Item(132, 11).Value = 607597.8  # K132
$ws.Cells.Item(132, 13).Value = -605067.8  # M132
$ws.Cells.Item(134, 8).Value = 3000  # H134
$ws.Cells.Item(134, 9).Value = 3000  # I134
$ws.Cells.Item(134, 10).Value = 3000  # J134
$ws.Cells.Item(134, 11).Value = 9000  # K134
$ws.Cells.Item(134, 12).Value = 9000  # L134
$ws.Cells.Item(134, 13).Value = -6465  # M134
$ws.Cells.Item(134, 14).Value = -14070  # N134
$ws.Cells.Item(136, 8).Value = 4226.857  # H136
$ws.Cells.Item(136, 9).Value = 4098  # I136
$ws.Cells.Item(136, 10).Value = 5000  # J136
$ws.Cells.Item(136, 11).Value = 12294  # K136
$ws.Cells.Item(136, 12).Value = 15000  # L136
$ws.Cells.Item(136, 13).Value = -9744  # M136
$ws.Cells.Item(136, 14).Value = -20100  # N136

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 0  # H69
$ws.Cells.Item(69, 9).Value = 0  # I69
$ws.Cells.Item(69, 11).Value = 0  # K69
$ws.Cells.Item(69, 13).ClearContents()  # M69
$ws.Cells.Item(72, 8).Value = 0  # H72
$ws.Cells.Item(72, 9).Value = 0  # I72
$ws.Cells.Item(72, 11).Value = 0  # K72
$ws.Cells.Item(72, 13).ClearContents()  # M72
$ws.Cells.Item(92, 8).Value = 333.33334  # H92
$ws.Cells.Item(92, 9).Value = 0  # I92
$ws.Cells.Item(92, 11).Value = 0  # K92
$ws.Cells.Item(92, 13).ClearContents()  # M92
$ws.Cells.Item(113, 8).Value = 1907.1818  # H113
$ws.Cells.Item(113, 9).Value = 995.2  # I113
$ws.Cells.Item(113, 10).Value = 2667.1667  # J113
$ws.Cells.Item(113, 11).Value = 2985.6  # K113
$ws.Cells.Item(113, 12).Value = 8001.500100000001  # L113
$ws.Cells.Item(113, 13).Value = -815.6000000000004  # M113
$ws.Cells.Item(113, 14).Value = -12341.5001  # N113

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 89999  # H64
$ws.Cells.Item(64, 10).Value = 89999  # J64
$ws.Cells.Item(64, 12).Value = 89999  # L64
$ws.Cells.Item(64, 14).Value = -90495  # N64
$ws.Cells.Item(67, 8).Value = 89999  # H67
$ws.Cells.Item(67, 10).Value = 89999  # J67
$ws.Cells.Item(67, 12).Value = 89999  # L67
$ws.Cells.Item(67, 14).Value = -91715  # N67
$ws.Cells.Item(70, 8).Value = 9250  # H70
$ws.Cells.Item(73, 8).Value = 9250  # H73
$ws.Cells.Item(122, 8).Value = 8721.4  # H122
$ws.Cells.Item(122, 9).Value = 5269  # I122
$ws.Cells.Item(122, 10).Value = 13900  # J122
$ws.Cells.Item(122, 11).Value = 15807  # K122
$ws.Cells.Item(122, 12).Value = 41700  # L122
$ws.Cells.Item(122, 13).Value = -13357  # M122
$ws.Cells.Item(122, 14).Value = -46600  # N122
$ws.Cells.Item(126, 8).Value = 4000  # H126
$ws.Cells.Item(126, 9).Value = 4000  # I126
$ws.Cells.Item(126, 11).Value = 12000  # K126
$ws.Cells.Item(126, 13).Value = -9530  # M126
$ws.Cells.Item(139, 8).Value = 0  # H139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 12).Value = 0  # L139
$ws.Cells.Item(139, 13).ClearContents()  # M139

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2911  # H7
$ws.Cells.Item(7, 9).Value = 1200  # I7
$ws.Cells.Item(7, 11).Value = 1200  # K7
$ws.Cells.Item(7, 13).Value = -1088  # M7
$ws.Cells.Item(40, 8).Value = 2974.8333  # H40
$ws.Cells.Item(40, 9).Value = 2968.8  # I40
$ws.Cells.Item(40, 11).Value = 2968.8  # K40
$ws.Cells.Item(40, 13).Value = -2832.8  # M40
$ws.Cells.Item(122, 8).Value = 1695  # H122
$ws.Cells.Item(122, 9).Value = 1596.6666  # I122
$ws.Cells.Item(122, 11).Value = 4789.9998  # K122
$ws.Cells.Item(122, 13).Value = -2339.9998  # M122
$ws.Cells.Item(126, 8).Value = 2911  # H126
$ws.Cells.Item(126, 9).Value = 1200  # I126
$ws.Cells.Item(126, 11).Value = 3600  # K126
$ws.Cells.Item(126, 13).Value = -1130  # M126
$ws.Cells.Item(132, 8).Value = 3915  # H132
$ws.Cells.Item(132, 9).Value = 2862.125  # I132
$ws.Cells.Item(132, 11).Value = 8586.375  # K132
$ws.Cells.Item(132, 13).Value = -6056.375  # M132
$ws.Cells.Item(136, 8).Value = 0  # H136
$ws.Cells.Item(136, 9).Value = 0  # I136
$ws.Cells.Item(136, 10).Value = 0  # J136
$ws.Cells.Item(136, 11).Value = 0  # K136
$ws.Cells.Item(136, 12).Value = 0  # L136
$ws.Cells.Item(136, 13).ClearContents()  # M136
$ws.Cells.Item(136, 14).ClearContents()  # N136

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 33702.832  # H74
$ws.Cells.Item(74, 10).Value = 31443.6  # J74
$ws.Cells.Item(74, 12).Value = 31443.6  # L74
$ws.Cells.Item(74, 14).Value = -33315.6  # N74
$ws.Cells.Item(77, 8).Value = 33702.832  # H77
$ws.Cells.Item(77, 10).Value = 31443.6  # J77
$ws.Cells.Item(77, 12).Value = 94330.79999999999  # L77
$ws.Cells.Item(77, 14).Value = -103690.8  # N77
$ws.Cells.Item(122, 8).Value = 3467.6667  # H122
$ws.Cells.Item(122, 9).Value = 2868.1667  # I122
$ws.Cells.Item(122, 10).Value = 4666.6665  # J122
$ws.Cells.Item(122, 11).Value = 8604.500100000001  # K122
$ws.Cells.Item(122, 12).Value = 13999.9995  # L122
$ws.Cells.Item(122, 13).Value = -6154.500100000001  # M122
$ws.Cells.Item(122, 14).Value = -18899.9995  # N122
$ws.Cells.Item(132, 8).Value = 2613.2856  # H132
$ws.Cells.Item(132, 9).Value = 1991.1428  # I132
$ws.Cells.Item(132, 11).Value = 5973.428400000001  # K132
$ws.Cells.Item(132, 13).Value = -3443.428400000001  # M132
